$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: 2/9, 1 hr / worked with data ---
$ws.Range("A7:B7").Copy()
$ws.Range("A8:B8").PasteSpecial(-4122)
$ws.Range("B8").Value = "worked with data"
$ws.Range("A8").Value = "2/9, 1 hr"

# --- Row 9: 2/10, 1 hr / worked on shiny ---
$ws.Range("A7:B7").Copy()
$ws.Range("A9:B9").PasteSpecial(-4122)
$ws.Range("B9").Value = "worked on shiny"
$ws.Range("A9").Value = "2/10, 1 hr"

# --- Row 10: 2/11, 4 hours / worked on map, got shiny app working, added time slider, wrote notes for class, worked on cleaning github ---
$ws.Range("A7:B7").Copy()
$ws.Range("A10:B10").PasteSpecial(-4122)
$ws.Range("A10").Value = "2/11, 4 hours"
$ws.Range("B10").Value = "worked on map, got shiny app working, added time slider, wrote notes for class, worked on cleaning github"

# --- Row heights (matches wrapped-text auto height in the source workbook) ---
$ws.Rows(8).RowHeight = 41.4
$ws.Rows(9).RowHeight = 27.6
$ws.Rows(10).RowHeight = 220.8

$excel.CutCopyMode = $false

# --- View / selection state ---
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("M10").Select() | Out-Null
